# Update cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.289.55"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "2.272.12"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'308.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'97.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("D10").Value = "'34.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.46%  "
$ws.Range("D11").Value = "'0.0811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").Value = "'6.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "2.625.18"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "'14.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "2.283.06"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "'0.787"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").Value = "42.176.42"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'12.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("D20").Value = "0.0₃0906"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Value = "'236.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("D24").Value = "'2.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'23.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("D28").Value = "'37.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.43%  "
$ws.Range("D29").Value = "'9.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").Value = "'163.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").Value = "'17.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'0.0732"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.33%  "
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("E38").Value = "  -2.80%  "
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("D41").Value = "'4.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("E42").Value = "  -4.49%  "
$ws.Range("D43").Value = "1.945.83"
$ws.Range("E43").Value = "  -2.93%  "
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").Value = "'18.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("D48").Value = "'54.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").Value = "2.497.67"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").Value = "'91.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").Value = "'71.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.58%  "
